$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.64"
$ws.Range("E2").Value = "'-4.59%"
$ws.Range("D3").Value = "'39.24"
$ws.Range("E3").Value = "'-8.46%"
$ws.Range("D4").Value = "'5.065"
$ws.Range("E4").Value = "'-2.56%"
$ws.Range("D6").Value = "'4.236"
$ws.Range("E6").Value = "'-1.83%"
$ws.Range("D7").Value = "'1.599"
$ws.Range("E7").Value = "'-11.17%"
$ws.Range("D8").Value = "'0.9141"
$ws.Range("E8").Value = "'-3.79%"
$ws.Range("D9").Value = "'0.1031"
$ws.Range("E9").Value = "'-8.12%"
$ws.Range("D10").Value = "'0.1742"
$ws.Range("E10").Value = "'-6.76%"
$ws.Range("D11").Value = "'0.08899"
$ws.Range("E11").Value = "'-5.89%"
$ws.Range("D12").Value = "'0.04435"
$ws.Range("E12").Value = "'-3.75%"
$ws.Range("E13").Value = "'-0.45%"
$ws.Range("D14").Value = "'0.001270"
$ws.Range("E14").Value = "'-2.21%"
$ws.Range("D15").Value = "'0.005835"
$ws.Range("E15").Value = "'0.75%"
$ws.Range("E16").Value = "'2,419.31%"
$ws.Range("D17").Value = "'3.357"
$ws.Range("E17").Value = "'-0.03%"
$ws.Range("E18").Value = "'-5.24%"
$ws.Range("E19").Value = "'-1.48%"
$ws.Range("D20").Value = "'7.030"
$ws.Range("E20").Value = "'-5.59%"
$ws.Range("E21").Value = "'-2.80%"
$ws.Range("D22").Value = "'0.2759"
$ws.Range("E22").Value = "'8.45%"
$ws.Range("D23").Value = "'0.04143"
$ws.Range("E23").Value = "'0.14%"
$ws.Range("D24").Value = "'0.001206"
$ws.Range("E24").Value = "'-3.23%"
$ws.Range("D25").Value = "'0.004025"
$ws.Range("E25").Value = "'-6.03%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'8.41%"
$ws.Range("D38").Value = "'0.02380"
$ws.Range("E38").Value = "'-10.63%"
$ws.Range("D39").Value = "'0.05162"
$ws.Range("E39").Value = "'-7.10%"
$ws.Range("D40").Value = "'0.007902"
$ws.Range("E40").Value = "'-2.77%"
$ws.Range("D41").Value = "'0.1315"
$ws.Range("E41").Value = "'-6.17%"
$ws.Range("E42").Value = "'-10.78%"
$ws.Range("E43").Value = "'-0.52%"
$ws.Range("D44").Value = "'0.007407"
$ws.Range("E44").Value = "'-3.00%"
$ws.Range("D45").Value = "'0.3327"
$ws.Range("E45").Value = "'3.65%"
$ws.Range("D46").Value = "'0.00006442"
$ws.Range("E46").Value = "'-4.44%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.31%"
$ws.Range("D48").Value = "'0.003003"
$ws.Range("E48").Value = "'-26.61%"
$ws.Range("D49").Value = "'0.004168"
$ws.Range("E49").Value = "'25.20%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.31%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.31%"

# Clear the quote-prefix formatting flag that typing a leading apostrophe
# sets, so the cells keep their original (default) style/number format
# while retaining the text values assigned above.
$ws.Range("D2:E51").Style = "Normal"
